# Swap the contents of columns D<->G and E<->F (codeforiati:group-name/
# codeforiati:category-name/codeforiati:group-code/codeforiati:category-code
# columns) for every row, including the header, to reorder the columns as:
#   code, name, status, codeforiati:category-code, codeforiati:group-code,
#   codeforiati:group-name, codeforiati:category-name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2.ToString()
    $eVal = $ws.Cells.Item($r, 5).Value2.ToString()
    $fVal = $ws.Cells.Item($r, 6).Value2.ToString()
    $gVal = $ws.Cells.Item($r, 7).Value2.ToString()

    # Re-assign with a leading apostrophe so Excel stores these as text
    # (matching the original shared-string / text cell type) rather than
    # auto-converting digit-only codes like "110" into numbers.
    $ws.Cells.Item($r, 4).Value2 = "'" + $gVal
    $ws.Cells.Item($r, 5).Value2 = "'" + $fVal
    $ws.Cells.Item($r, 6).Value2 = "'" + $dVal
    $ws.Cells.Item($r, 7).Value2 = "'" + $eVal
}
